# Update timestamps in the handback-status report, as part of
# regenerating the report for handback.

$wb = $excel.ActiveWorkbook

# Overview sheet: "Latest HO Xliff Generate Date" for the
# 3074847d-c8c4-4e29-b8c5-dc8af17a4364.md row
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G4").Value = "2016-08-27 04:43:57"

# zh-cn sheet: "Correspond Handoff Datetime" / "Correspond Handback DateTime"
# for the 3074847d-c8c4-4e29-b8c5-dc8af17a4364 row
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H4").Value = "2016-08-27 04:43:53"
$wsZhCn.Range("K4").Value = "2016-08-27 04:44:14"

# de-de sheet: "Correspond Handback DateTime" for the
# 3074847d-c8c4-4e29-b8c5-dc8af17a4364 row
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("K4").Value = "2016-08-27 04:44:21"
